# Update the NATMI TPM output values (columns G:T, rows 2-11) with newly
# recomputed TPM-based statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @{
    2  = @{ G=31.35623066666667;  H=94.068692;  I=0.07215642027787079; J=0.07299614919666826; K=1; L=0.3333333333333333; M=0.1375686666666667; N=0.412706; O=0.2896572731203081; P=0.2896572731203081; Q=4.313634844505778;  R=38.822713600552;  S=0.02090063193581096;  T=0.0211438655245901 }
    3  = @{ G=31.35623066666667;  H=94.068692;  I=0.07215642027787079; J=0.07299614919666826; K=2; L=0.6666666666666666; M=0.3373673333333334; N=1.012102; O=0.7103427268796919; P=0.7103427268796919; Q=10.57856792339822;  R=95.207111310584;  S=0.05125578834205984;  T=0.05185228367207816 }
    4  = @{                                                            I=0.3221618650682612;  J=0.3259110620709639;  K=1; L=0.3333333333333333; M=0.1375686666666667; N=0.412706; O=0.2896572731203081; P=0.2896572731203081; Q=19.25939010524356;  R=173.334510947192; S=0.09331652733902521; T=0.09440250951921889 }
    5  = @{                                                            I=0.3221618650682612;  J=0.3259110620709639;  K=2; L=0.6666666666666666; M=0.3373673333333334; N=1.012102; O=0.7103427268796919; P=0.7103427268796919; Q=47.23087923194045;  R=425.077913087464; S=0.228845337729236;   T=0.231508552551745 }
    6  = @{ G=130.001713;         H=390.005139; I=0.2991577125385502;  J=0.3026391959814998;  K=1; L=0.3333333333333333; M=0.1375686666666667; N=0.412706; O=0.2896572731203081; P=0.2896572731203081; Q=17.88416232179267;  R=160.957460896134; S=0.08665320724682547; T=0.08766164424732376 }
    7  = @{ G=130.001713;         H=390.005139; I=0.2991577125385502;  J=0.3026391959814998;  K=2; L=0.6666666666666666; M=0.3373673333333334; N=1.012102; O=0.7103427268796919; P=0.7103427268796919; Q=43.85833124357534;  R=394.724981192178; S=0.2125045052917248;  T=0.2149775517341761 }
    8  = @{ G=14.997169;          H=29.994338;  I=0.03451122811430998; J=0.02327523775607825; K=1; L=0.3333333333333333; M=0.1375686666666667; N=0.412706; O=0.2896572731203081; P=0.2896572731203081; Q=2.063140543104667;  R=12.378843258628;  S=0.009996428227623942; T=0.006741841899652466 }
    9  = @{ G=14.997169;          H=29.994338;  I=0.03451122811430998; J=0.02327523775607825; K=2; L=0.6666666666666666; M=0.3373673333333334; N=1.012102; O=0.7103427268796919; P=0.7103427268796919; Q=5.059554913079333;  R=30.357329478476;  S=0.02451479988668604;  T=0.01653339585642579 }
    10 = @{ G=118.2056323333333;  H=354.616897; I=0.2720127740010079;  J=0.2751783549947898;  K=1; L=0.3333333333333333; M=0.1375686666666667; N=0.412706; O=0.2896572731203081; P=0.2896572731203081; Q=16.26139123258689;  R=146.352521093282; S=0.0787904783710226;  T=0.07970741192952294 }
    11 = @{ G=118.2056323333333;  H=354.616897; I=0.2720127740010079;  J=0.2751783549947898;  K=2; L=0.6666666666666666; M=0.3373673333333334; N=1.012102; O=0.7103427268796919; P=0.7103427268796919; Q=39.87871896527711;  R=358.908470687494; S=0.1932222956299853;  T=0.1954709430652669 }
}

foreach ($rowNum in $newData.Keys) {
    $rowVals = $newData[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
